$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{"A"="DANATEST"; "B"=2; "C"=1; "D"="Best"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=43; "J"=36; "K"=32; "L"=124; "M"=126; "N"=0.64481394162585348; "O"=0.01; "P"=10; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=1; "D"="MaxLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=43; "J"=39; "K"=34; "L"=114; "M"=126; "N"=0.67481394162585351; "O"=0.005; "P"=10; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=1; "D"="MinLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=43; "J"=48; "K"=29; "L"=145; "M"=126; "N"=0.58481394162585343; "O"=0.005; "P"=10; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=2; "D"="Best"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=44; "J"=16; "K"=31; "L"=134; "M"=135; "N"=0.61811682663654688; "O"=0.005; "P"=10; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=2; "D"="MaxLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=44; "J"=21; "K"=33; "L"=121; "M"=135; "N"=0.65311682663654691; "O"=0.005; "P"=10; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=2; "D"="MinLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=44; "J"=25; "K"=30; "L"=139; "M"=135; "N"=0.60311682663654687; "O"=0.005; "P"=10; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=3; "D"="Best"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=44; "J"=39; "K"=32; "L"=128; "M"=118; "N"=0.63409770834458379; "O"=0.005; "P"=15; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=3; "D"="MaxLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=44; "J"=45; "K"=33; "L"=119; "M"=118; "N"=0.65909770834458381; "O"=0.005; "P"=15; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=3; "D"="MinLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=44; "J"=49; "K"=31; "L"=135; "M"=118; "N"=0.61409770834458377; "O"=0.005; "P"=15; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=4; "D"="Best"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=45; "J"=10; "K"=32; "L"=129; "M"=124; "N"=0.6303408694257544; "O"=0.005; "P"=1; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=4; "D"="MaxLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=45; "J"=13; "K"=33; "L"=122; "M"=124; "N"=0.65034086942575442; "O"=0.005; "P"=1; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=4; "D"="MinLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=45; "J"=20; "K"=31; "L"=136; "M"=124; "N"=0.61034086942575438; "O"=0.005; "P"=1; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=5; "D"="Best"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=45; "J"=52; "K"=31; "L"=130; "M"=131; "N"=0.62743452371495156; "O"=0.005; "P"=1; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=5; "D"="MaxLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=45; "J"=56; "K"=33; "L"=120; "M"=131; "N"=0.65743452371495159; "O"=0.005; "P"=1; "Q"=1},
    @{"A"="DANATEST"; "B"=2; "C"=5; "D"="MinLambda"; "E"=2023; "F"=11; "G"=22; "H"=14; "I"=46; "J"=0; "K"=30; "L"=137; "M"=131; "N"=0.60743452371495155; "O"=0.005; "P"=1; "Q"=1},
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    foreach ($col in $cols) {
        $ws.Range($col + $r).Value = $rowData[$col]
    }
}
